# Auto-generated edit script applying scheduled market-data refresh to the
# per-job Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Each sheet carries a native Excel Table (Table_<JOB>) over A1:N141 with
# columns: ... H=currentAveragePrice, I=currentAveragePriceNQ,
# J=currentAveragePriceHQ, K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ,
# N=LeveProfitHQ. Values below are freshly recomputed snapshot values (no
# formulas are stored in this workbook), so we just overwrite H:N per row.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 666.3333
$ws.Range("I29").Value = 666.3333
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 1998.9999
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -1717.9999
$ws.Range("N29").ClearContents()
$ws.Range("H70").Value = 1333.9286
$ws.Range("I70").Value = 993.5484
$ws.Range("J70").Value = 3971.875
$ws.Range("K70").Value = 2980.6452
$ws.Range("L70").Value = 11915.625
$ws.Range("M70").Value = -2710.6452
$ws.Range("N70").Value = -12455.625
$ws.Range("H73").Value = 1333.9286
$ws.Range("I73").Value = 993.5484
$ws.Range("J73").Value = 3971.875
$ws.Range("K73").Value = 2980.6452
$ws.Range("L73").Value = 11915.625
$ws.Range("M73").Value = -2044.6452
$ws.Range("N73").Value = -13787.625
$ws.Range("H98").Value = 4007.0417
$ws.Range("I98").Value = 2429.375
$ws.Range("J98").Value = 7162.375
$ws.Range("K98").Value = 2429.375
$ws.Range("L98").Value = 7162.375
$ws.Range("M98").Value = -931.375
$ws.Range("N98").Value = -10158.375
$ws.Range("H122").Value = 4007.0417
$ws.Range("I122").Value = 2429.375
$ws.Range("J122").Value = 7162.375
$ws.Range("K122").Value = 7288.125
$ws.Range("L122").Value = 21487.125
$ws.Range("M122").Value = -4838.125
$ws.Range("N122").Value = -26387.125
$ws.Range("H137").Value = 3541.4348
$ws.Range("I137").Value = 2331.077
$ws.Range("J137").Value = 5114.9
$ws.Range("K137").Value = 6993.231000000001
$ws.Range("L137").Value = 15344.7
$ws.Range("M137").Value = -4443.231000000001
$ws.Range("N137").Value = -20444.7
$ws.Range("H138").Value = 2053.1272
$ws.Range("I138").Value = 675.1818
$ws.Range("J138").Value = 2971.7576
$ws.Range("K138").Value = 2025.5454
$ws.Range("L138").Value = 8915.272799999999
$ws.Range("M138").Value = 3114.4546
$ws.Range("N138").Value = -19195.2728

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7110.2593
$ws.Range("I32").Value = 5276.029
$ws.Range("J32").Value = 10489.105
$ws.Range("K32").Value = 5276.029
$ws.Range("L32").Value = 10489.105
$ws.Range("M32").Value = -4989.029
$ws.Range("N32").Value = -11063.105
$ws.Range("H122").Value = 2336.6086
$ws.Range("I122").Value = 1331.8334
$ws.Range("J122").Value = 3432.7273
$ws.Range("K122").Value = 3995.5002
$ws.Range("L122").Value = 10298.1819
$ws.Range("M122").Value = -1545.5002
$ws.Range("N122").Value = -15198.1819
$ws.Range("H125").Value = 43354
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 43354
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 43354
$ws.Range("N125").Value = -53194
$ws.Range("H137").Value = 40577.5
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 40577.5
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 40577.5
$ws.Range("N137").Value = -50777.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 118880
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 118880
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 118880
$ws.Range("N59").Value = -120574
$ws.Range("H99").Value = 2101.8462
$ws.Range("I99").Value = 1526.8125
$ws.Range("J99").Value = 3021.9
$ws.Range("K99").Value = 1526.8125
$ws.Range("L99").Value = 3021.9
$ws.Range("M99").Value = -28.8125
$ws.Range("N99").Value = -6017.9
$ws.Range("H134").Value = 2502.9778
$ws.Range("I134").Value = 1311.079
$ws.Range("J134").Value = 8973.286
$ws.Range("K134").Value = 3933.237
$ws.Range("L134").Value = 26919.858
$ws.Range("M134").Value = -1398.237
$ws.Range("N134").Value = -31989.858
$ws.Range("H137").Value = 38784.445
$ws.Range("I137").Value = 25000
$ws.Range("J137").Value = 40507.5
$ws.Range("K137").Value = 25000
$ws.Range("L137").Value = 40507.5
$ws.Range("M137").Value = -19900
$ws.Range("N137").Value = -50707.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14709604
$ws.Range("I31").Value = 1994.2667
$ws.Range("J31").Value = 26320874
$ws.Range("K31").Value = 1994.2667
$ws.Range("L31").Value = 26320874
$ws.Range("M31").Value = -1699.2667
$ws.Range("N31").Value = -26321464
$ws.Range("H33").Value = 8499.5
$ws.Range("I33").Value = 8499.5
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 8499.5
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -8120.5
$ws.Range("H34").Value = 14709604
$ws.Range("I34").Value = 1994.2667
$ws.Range("J34").Value = 26320874
$ws.Range("K34").Value = 1994.2667
$ws.Range("L34").Value = 26320874
$ws.Range("M34").Value = -1792.2667
$ws.Range("N34").Value = -26321278
$ws.Range("H62").Value = 4950
$ws.Range("I62").Value = 5800
$ws.Range("J62").Value = 2400
$ws.Range("K62").Value = 5800
$ws.Range("L62").Value = 2400
$ws.Range("M62").Value = -5176
$ws.Range("N62").Value = -3648
$ws.Range("H65").Value = 4950
$ws.Range("I65").Value = 5800
$ws.Range("J65").Value = 2400
$ws.Range("K65").Value = 29000
$ws.Range("L65").Value = 12000
$ws.Range("M65").Value = -25880
$ws.Range("N65").Value = -18240
$ws.Range("H99").Value = 16672392
$ws.Range("I99").Value = 28574816
$ws.Range("J99").Value = 8998
$ws.Range("K99").Value = 28574816
$ws.Range("L99").Value = 8998
$ws.Range("M99").Value = -28573318
$ws.Range("H122").Value = 2563
$ws.Range("I122").Value = 1989.2307
$ws.Range("J122").Value = 3628.5715
$ws.Range("K122").Value = 5967.6921
$ws.Range("L122").Value = 10885.7145
$ws.Range("M122").Value = -3517.6921
$ws.Range("H126").Value = 16672392
$ws.Range("I126").Value = 28574816
$ws.Range("J126").Value = 8998
$ws.Range("K126").Value = 85724448
$ws.Range("L126").Value = 26994
$ws.Range("M126").Value = -85721978

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 569.58826
$ws.Range("I113").Value = 558.36
$ws.Range("J113").Value = 580.38464
$ws.Range("K113").Value = 1675.08
$ws.Range("L113").Value = 1741.15392
$ws.Range("M113").Value = 494.9200000000001
$ws.Range("N113").Value = -6081.15392
$ws.Range("H125").Value = 7666.6665
$ws.Range("I125").Value = 3000
$ws.Range("J125").Value = 10000
$ws.Range("K125").Value = 9000
$ws.Range("L125").Value = 30000
$ws.Range("M125").Value = -4080
$ws.Range("N125").Value = -39840
$ws.Range("H139").Value = 1488.75
$ws.Range("I139").Value = 1321.3334
$ws.Range("J139").Value = 4000
$ws.Range("K139").Value = 3964.0002
$ws.Range("L139").Value = 12000
$ws.Range("M139").Value = 1175.9998
$ws.Range("N139").Value = -22280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 33549.715
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 33549.715
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 33549.715
$ws.Range("N64").Value = -34045.715
$ws.Range("H67").Value = 33549.715
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 33549.715
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 33549.715
$ws.Range("N67").Value = -35265.715
$ws.Range("H97").Value = 699.9167
$ws.Range("I97").Value = 735.2727
$ws.Range("J97").Value = 311
$ws.Range("K97").Value = 735.2727
$ws.Range("L97").Value = 311
$ws.Range("M97").Value = -239.2727
$ws.Range("N97").Value = -1303
$ws.Range("H102").Value = 2094.3713
$ws.Range("I102").Value = 1668.2963
$ws.Range("J102").Value = 3532.375
$ws.Range("K102").Value = 1668.2963
$ws.Range("L102").Value = 3532.375
$ws.Range("M102").Value = -46.29629999999997
$ws.Range("H119").Value = 40400
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 40400
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 40400
$ws.Range("N119").Value = -50076
$ws.Range("H122").Value = 3228.7856
$ws.Range("I122").Value = 2012.4375
$ws.Range("J122").Value = 4850.5835
$ws.Range("K122").Value = 6037.3125
$ws.Range("L122").Value = 14551.7505
$ws.Range("M122").Value = -3587.3125
$ws.Range("N122").Value = -19451.7505

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7955.8945
$ws.Range("I40").Value = 5996.2
$ws.Range("J40").Value = 10133.333
$ws.Range("K40").Value = 5996.2
$ws.Range("L40").Value = 10133.333
$ws.Range("M40").Value = -5860.2
$ws.Range("N40").Value = -10405.333
$ws.Range("H100").Value = 1984
$ws.Range("I100").Value = 1799.2727
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 1799.2727
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -1258.2727
$ws.Range("N100").Value = -4082
$ws.Range("H122").Value = 4038.3215
$ws.Range("I122").Value = 2242.0476
$ws.Range("J122").Value = 9427.143
$ws.Range("K122").Value = 6726.1428
$ws.Range("L122").Value = 28281.429
$ws.Range("M122").Value = -4276.1428
$ws.Range("N122").Value = -33181.429
$ws.Range("H132").Value = 4899.724
$ws.Range("I132").Value = 1884.1154
$ws.Range("J132").Value = 7349.9062
$ws.Range("K132").Value = 5652.3462
$ws.Range("L132").Value = 22049.7186
$ws.Range("M132").Value = -3122.3462
$ws.Range("N132").Value = -27109.7186

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("N78").ClearContents()
$ws.Range("H122").Value = 2446.175
$ws.Range("I122").Value = 1339.5518
$ws.Range("J122").Value = 5363.636
$ws.Range("K122").Value = 4018.6554
$ws.Range("L122").Value = 16090.908
$ws.Range("M122").Value = -1568.6554
$ws.Range("N122").Value = -20990.908
$ws.Range("H126").Value = 2800.2593
$ws.Range("I126").Value = 1912.4615
$ws.Range("J126").Value = 3624.6428
$ws.Range("K126").Value = 5737.3845
$ws.Range("L126").Value = 10873.9284
$ws.Range("M126").Value = -3267.3845
$ws.Range("H136").Value = 1394.2424
$ws.Range("I136").Value = 883.5263
$ws.Range("J136").Value = 2087.3572
$ws.Range("K136").Value = 2650.5789
$ws.Range("L136").Value = 6262.071599999999
$ws.Range("M136").Value = -100.5789
$ws.Range("N136").Value = -11362.0716

